$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("DOLCYL 2MG 30 TAB."): current-balance and transaction-count text
# values refreshed for the new report pull.
$ws.Range("H13").Value = "0:1"
$ws.Range("Q13").Value = "0:3"

# Selling price (col P, stored as text) updated to the new figure. Flip the
# cell to text format first so Excel keeps storing it as a string instead of
# re-typing it as a number, then restore the original number format so the
# cell's style/formatting stays untouched.
$fmt = $ws.Range("P13").NumberFormat
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "47.5200"
$ws.Range("P13").NumberFormat = $fmt

# Selling-price column total (P20) reflects the updated P13 figure.
$ws.Range("P20").Value = 792.05

# Footer timestamp bumped to the new export time.
$ws.Range("A21").Value = "Sunday, 8 June, 2025 10:47 AM"
